$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a price into the D column while forcing it to stay a text
# value (many of these look like plain decimal numbers, e.g. "695.63",
# and Excel would otherwise silently convert them to numeric cells).
function Set-PriceText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-PriceText "D2" "71.164.48"
$ws.Range("E2").Value = "  +0.83%  "

# Row 3 - Ethereum
Set-PriceText "D3" "3.847.42"
$ws.Range("E3").Value = "  +0.81%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-PriceText "D5" "695.63"
$ws.Range("E5").Value = "  +1.40%  "

# Row 6 - Solana
Set-PriceText "D6" "173.65"
$ws.Range("E6").Value = "  +2.10%  "

# Row 7 - LidoStakedEther
Set-PriceText "D7" "3.845.23"
$ws.Range("E7").Value = "  +0.78%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.10%  "

# Row 11 - Toncoin
Set-PriceText "D11" "7.27"
$ws.Range("E11").Value = "  +1.06%  "

# Row 12 - Cardano
Set-PriceText "D12" "0.461"
$ws.Range("E12").Value = "  -0.16%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +4.78%  "

# Row 14 - Avalanche
Set-PriceText "D14" "36.42"
$ws.Range("E14").Value = "  +1.58%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-PriceText "D15" "4.496.38"
$ws.Range("E15").Value = "  +0.83%  "

# Row 16 - WrappedEther
Set-PriceText "D16" "3.849.09"
$ws.Range("E16").Value = "  +0.83%  "

# Row 17 - WrappedBTC
Set-PriceText "D17" "71.234.22"
$ws.Range("E17").Value = "  +0.79%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  +0.05%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +0.49%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +0.13%  "

# Row 21 - Uniswap
Set-PriceText "D21" "11.18"
$ws.Range("E21").Value = "  -0.58%  "

# Row 22 - BitcoinCash
Set-PriceText "D22" "493.12"

# Row 23 - Polygon
$ws.Range("E23").Value = "  +1.19%  "

# Row 24 - Litecoin
Set-PriceText "D24" "85.05"
$ws.Range("E24").Value = "  +2.10%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +1.91%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("E26").Value = "  +0.27%  "

# Row 27 - RenderToken
$ws.Range("E27").Value = "  +2.49%  "

# Row 28 - Fetch.AI
$ws.Range("E28").Value = "  +1.84%  "

# Row 29 - WrappedeETH
Set-PriceText "D29" "4.002.95"
$ws.Range("E29").Value = "  +0.88%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +8.13%  "

# Row 31 - Dai
$ws.Range("E31").Value = "  -0.13%  "

# Row 32 - NEARProtocol
Set-PriceText "D32" "7.62"
$ws.Range("E32").Value = "  +2.82%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  -0.33%  "

# Row 34 - EthereumClassic
Set-PriceText "D34" "29.65"
$ws.Range("E34").Value = "  +0.24%  "

# Row 35 - Kaspa
Set-PriceText "D35" "0.181"
$ws.Range("E35").Value = "  +0.12%  "

# Row 36 - Aptos
$ws.Range("E36").Value = "  +1.43%  "

# Row 37 - RenzoRestakedETH
Set-PriceText "D37" "3.799.55"
$ws.Range("E37").Value = "  +0.69%  "

# Row 38 - Binance-PegBSC-USD
Set-PriceText "D38" "0.999"
$ws.Range("E38").Value = "  -0.09%  "

# Row 39 - Hedera
$ws.Range("E39").Value = "  +2.59%  "

# Row 40 - Stacks
$ws.Range("E40").Value = "  +12.56%  "

# Row 41 - dogwifhat
$ws.Range("E41").Value = "  +0.21%  "

# Row 42 - Filecoin
Set-PriceText "D42" "6.05"
$ws.Range("E42").Value = "  +1.72%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  +6.42%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  -0.05%  "

# Row 46 - Monero
Set-PriceText "D46" "163.33"
$ws.Range("E46").Value = "  +2.39%  "

# Row 47 - FLOKI
$ws.Range("E47").Value = "  +2.57%  "

# Row 48 - OKB
$ws.Range("E48").Value = "  +0.98%  "

# Row 49 - Arweave
Set-PriceText "D49" "44.24"
$ws.Range("E49").Value = "  -3.85%  "

# Row 50 - was Bittensor, now TheGraph
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-PriceText "D50" "0.303"
$ws.Range("E50").Value = "  +1.12%  "

# Row 51 - was TheGraph, now Bittensor
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-PriceText "D51" "418.75"
$ws.Range("E51").Value = "  +5.65%  "
